# Generate output files and plots
# Update cap-factor style input data on "elec_demand (1)" and "elec_demand (2)"
# and refresh view/selection state to match the authored commit.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("rent_cap")
$ws8 = $wb.Worksheets.Item("elec_demand (1)")
$ws9 = $wb.Worksheets.Item("elec_demand (2)")

# --- elec_demand (1): overwrite rows 2-4 with newly generated values ---
    $ws8.Range("B2").Value = 0.5223607840973403
    $ws8.Range("C2").Value = 0.4976949585403172
    $ws8.Range("D2").Value = 0.4420770719465928
    $ws8.Range("E2").Value = 0.4377774266944169
    $ws8.Range("F2").Value = 0.4352162470552045
    $ws8.Range("G2").Value = 0.44180383833619113
    $ws8.Range("H2").Value = 0.5410607685473495
    $ws8.Range("I2").Value = 0.42634358331726746
    $ws8.Range("J2").Value = 0.05518583660644147
    $ws8.Range("K2").Value = 0
    $ws8.Range("L2").Value = 0
    $ws8.Range("M2").Value = 0.8571057128225469
    $ws8.Range("N2").Value = 0.7515731623147012
    $ws8.Range("O2").Value = 0.1106265662584133
    $ws8.Range("P2").Value = 0.00045206305619120257
    $ws8.Range("Q2").Value = 0.2721617076023392
    $ws8.Range("R2").Value = 0.8589286852955188
    $ws8.Range("S2").Value = 0.777598571259276
    $ws8.Range("T2").Value = 0.6923943458764087
    $ws8.Range("U2").Value = 0.7064915990693572
    $ws8.Range("V2").Value = 0.6593004021253852
    $ws8.Range("W2").Value = 0.6837411441237511
    $ws8.Range("X2").Value = 0.6465764204867005
    $ws8.Range("Y2").Value = 0.5278428793413089
    $ws8.Range("B3").Value = 0.4532206627457997
    $ws8.Range("C3").Value = 0.42851577557190434
    $ws8.Range("D3").Value = 0.410062670020118
    $ws8.Range("E3").Value = 0.4049464029192532
    $ws8.Range("F3").Value = 0.3978478327247562
    $ws8.Range("G3").Value = 0.3806088878530444
    $ws8.Range("H3").Value = 0.4409417020283103
    $ws8.Range("I3").Value = 0.3698487681635868
    $ws8.Range("J3").Value = 0.1610332845513586
    $ws8.Range("K3").Value = 0.16500218302712158
    $ws8.Range("L3").Value = 0.1795243482064742
    $ws8.Range("M3").Value = 0.605897288607249
    $ws8.Range("N3").Value = 0.5524125623034574
    $ws8.Range("O3").Value = 0.2224012234616651
    $ws8.Range("P3").Value = 0.1661176493321201
    $ws8.Range("Q3").Value = 0.2997708038330913
    $ws8.Range("R3").Value = 0.5992202640346023
    $ws8.Range("S3").Value = 0.5283016559623788
    $ws8.Range("T3").Value = 0.4354426497707208
    $ws8.Range("U3").Value = 0.4004820172299548
    $ws8.Range("V3").Value = 0.42149519720844475
    $ws8.Range("W3").Value = 0.4546003728097659
    $ws8.Range("X3").Value = 0.4526678540640823
    $ws8.Range("Y3").Value = 0.4251154849589451
    $ws8.Range("B4").Value = 0.3840805413942591
    $ws8.Range("C4").Value = 0.3593365926034915
    $ws8.Range("D4").Value = 0.3780482680936434
    $ws8.Range("E4").Value = 0.3721153791440894
    $ws8.Range("F4").Value = 0.360479418394308
    $ws8.Range("G4").Value = 0.3194139373698977
    $ws8.Range("H4").Value = 0.340822635509271
    $ws8.Range("I4").Value = 0.313353953009906
    $ws8.Range("J4").Value = 0.2668807324962758
    $ws8.Range("K4").Value = 0.33000436605424316
    $ws8.Range("L4").Value = 0.3590486964129484
    $ws8.Range("M4").Value = 0.35468886439195096
    $ws8.Range("N4").Value = 0.35325196229221345
    $ws8.Range("O4").Value = 0.33417588066491694
    $ws8.Range("P4").Value = 0.331783235608049
    $ws8.Range("Q4").Value = 0.32737990006384343
    $ws8.Range("R4").Value = 0.3395118427736855
    $ws8.Range("S4").Value = 0.27900474066548137
    $ws8.Range("T4").Value = 0.178490953665033
    $ws8.Range("U4").Value = 0.0944724353905524
    $ws8.Range("V4").Value = 0.1836899922915041
    $ws8.Range("W4").Value = 0.2254596014957808
    $ws8.Range("X4").Value = 0.2587592876414642
    $ws8.Range("Y4").Value = 0.3223880905765812

# Rows 3 & 4 switch from the plain "0.00" style to the wrapped "0.00" style
# (same numeric format, matches row 2's look).
$ws8.Range("B3:Y4").WrapText = $true
$ws8.Range("B3:Y4").VerticalAlignment = -4108

# --- elec_demand (2): overwrite rows 2-4 with newly generated values ---
    $ws9.Range("B2").Value = 0.3415913978494623
    $ws9.Range("C2").Value = 0.3032688172043011
    $ws9.Range("D2").Value = 0.29459139784946226
    $ws9.Range("E2").Value = 0.28721505376344075
    $ws9.Range("F2").Value = 0.29210752688172037
    $ws9.Range("G2").Value = 0.2961075268817204
    $ws9.Range("H2").Value = 0.3111290322580645
    $ws9.Range("I2").Value = 0.345741935483871
    $ws9.Range("J2").Value = 0.3868387096774193
    $ws9.Range("K2").Value = 0.41622580645161295
    $ws9.Range("L2").Value = 0.45535483870967736
    $ws9.Range("M2").Value = 0.4138709677419354
    $ws9.Range("N2").Value = 0.3768494623655914
    $ws9.Range("O2").Value = 0.3618279569892473
    $ws9.Range("P2").Value = 0.3337956989247312
    $ws9.Range("Q2").Value = 0.3368387096774193
    $ws9.Range("R2").Value = 0.3406989247311828
    $ws9.Range("S2").Value = 0.3688709677419354
    $ws9.Range("T2").Value = 0.404494623655914
    $ws9.Range("U2").Value = 0.4571935483870968
    $ws9.Range("V2").Value = 0.479989247311828
    $ws9.Range("W2").Value = 0.4640860215053763
    $ws9.Range("X2").Value = 0.4129677419354838
    $ws9.Range("Y2").Value = 0.40031182795698916
    $ws9.Range("B3").Value = 0.2199123540971449
    $ws9.Range("C3").Value = 0.2022554430849092
    $ws9.Range("D3").Value = 0.1989724805339266
    $ws9.Range("E3").Value = 0.1943889349276974
    $ws9.Range("F3").Value = 0.1977337921764924
    $ws9.Range("G3").Value = 0.1996330737856878
    $ws9.Range("H3").Value = 0.201662763255469
    $ws9.Range("I3").Value = 0.2559364562476827
    $ws9.Range("J3").Value = 0.2959026019651465
    $ws9.Range("K3").Value = 0.716865661846496
    $ws9.Range("L3").Value = 0.9860919883203562
    $ws9.Range("M3").Value = 1.107966317204301
    $ws9.Range("N3").Value = 1.1344969150908422
    $ws9.Range("O3").Value = 1.1773252428624401
    $ws9.Range("P3").Value = 1.1448735965888028
    $ws9.Range("Q3").Value = 1.0590023433444569
    $ws9.Range("R3").Value = 0.9132089738598443
    $ws9.Range("S3").Value = 0.4577487035467642
    $ws9.Range("T3").Value = 0.32752397849462356
    $ws9.Range("U3").Value = 0.3321666714157706
    $ws9.Range("V3").Value = 0.33377549310035853
    $ws9.Range("W3").Value = 0.3265952413082437
    $ws9.Range("X3").Value = 0.300951482078853
    $ws9.Range("Y3").Value = 0.2941887723118279
    $ws9.Range("B4").Value = 0.09823331034482759
    $ws9.Range("C4").Value = 0.10124206896551721
    $ws9.Range("D4").Value = 0.1033535632183908
    $ws9.Range("E4").Value = 0.101562816091954
    $ws9.Range("F4").Value = 0.1033600574712644
    $ws9.Range("G4").Value = 0.1031586206896552
    $ws9.Range("H4").Value = 0.09219649425287356
    $ws9.Range("I4").Value = 0.16613097701149432
    $ws9.Range("J4").Value = 0.20496649425287358
    $ws9.Range("K4").Value = 1.017505517241379
    $ws9.Range("L4").Value = 1.5168291379310348
    $ws9.Range("M4").Value = 1.802061666666666
    $ws9.Range("N4").Value = 1.892144367816092
    $ws9.Range("O4").Value = 1.992822528735632
    $ws9.Range("P4").Value = 1.955951494252874
    $ws9.Range("Q4").Value = 1.781165977011494
    $ws9.Range("R4").Value = 1.485719022988506
    $ws9.Range("S4").Value = 0.546626439351593
    $ws9.Range("T4").Value = 0.2505533333333333
    $ws9.Range("U4").Value = 0.2071397944444444
    $ws9.Range("V4").Value = 0.18756173888888888
    $ws9.Range("W4").Value = 0.1891044611111111
    $ws9.Range("X4").Value = 0.1889352222222222
    $ws9.Range("Y4").Value = 0.18806571666666658

# The leftover placeholder cells on row 9 are no longer used.
$ws9.Range("B9:D9").Clear()

# --- view/selection state ---
$ws8.Range("B6").Select()

$ws9.Activate()
$ws9.Range("L12").Select()
